# "writeup updated; more results"
# Adds a third "50 moves" baseline run block (rows 17-19 headers + rows 20-39
# seed/score data + row 41 overall averages) to the "Multiple runs" sheet,
# plus a new "50 moves" label in B1.
#
# Cell writes below are deliberately ordered to match the original authoring
# sequence (first-use order of each new label), so new shared-string entries
# land in the same order as the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Multiple runs")

# --- Row 17: second "Baseline" sub-title + "200 moves" label ---
$ws.Range("A17").NumberFormat = "0"
$ws.Range("A17").Value = "Baseline "
$ws.Range("B17").NumberFormat = "0.00"
$ws.Range("B17").Value = "200 moves"

# --- Row 41 "overall" label (first use of the word, entered early) ---
$ws.Range("A41").NumberFormat = "0"
$ws.Range("A41").Value = "overall"

# --- Row 18: parameter labels for the three new result columns ---
$ws.Range("A18").NumberFormat = "0"
$ws.Range("A18").Value = "Thres=-100, scaredTime=10"
$ws.Range("D18").Value = "Thres=-125, time=10"

# --- Row 1: new "50 moves" label next to the existing "Baseline" title ---
$ws.Range("B1").NumberFormat = "0.00"
$ws.Range("B1").Value = "50 moves"

$ws.Range("G18").Value = "Thres=-150, time=10"

# --- Row 19: column headers (Seed / Avg Score) repeated per block ---
$ws.Range("A19").NumberFormat = "0"
$ws.Range("A19").Value = "Seed"
$ws.Range("B19").NumberFormat = "0.00"
$ws.Range("B19").Value = "Avg Score"

$ws.Range("D19").Value = "Seed"
$ws.Range("E19").NumberFormat = "0.00"
$ws.Range("E19").Value = "Avg Score"

$ws.Range("G19").Value = "Seed"
$ws.Range("H19").Value = "Avg Score"

# --- Rows 20-39: per-seed results for the first two blocks (A/B, D/E) ---
$data = @(
    @(20, 1, 1264.13553112, 1, 1264.13553112),
    @(21, 2, 1092.4423876599999, 2, 1092.4423876599999),
    @(22, 3, 648.60251877799999, 3, 648.60251877799999),
    @(23, 4, 255.41149254199999, 4, 255.41149254199999),
    @(24, 5, 30.087156974500001, 5, 579.61181571700001),
    @(25, 6, 999.81892902100003, 6, 999.81892902100003),
    @(26, 7, 862.66004696499999, 7, 562.50617603499995),
    @(27, 8, 600.16129394200004, 8, 600.16129394200004),
    @(28, 9, 901.04026481000005, 9, 901.04026481000005),
    @(29, 10, 467.28570589899999, 10, 818.49615234400005),
    @(30, 11, 1090.1498091399999, 11, 1086.03281528),
    @(31, 12, 1115.01770145, 12, 839.06133802199997),
    @(32, 13, 814.80119441399995, 13, 736.97497188399996),
    @(33, 14, 909.62100499500002, 14, 909.62100499500002),
    @(34, 15, 2049.4600154, 15, 2049.4600154),
    @(35, 16, 942.70378424900002, 16, 1084.4743467799999),
    @(36, 17, 531.50877935200003, 17, 845.61834408200002),
    @(37, 18, 313.75163984, 18, 313.75163984),
    @(38, 19, 1639.8191200199999, 19, 1639.8191200199999),
    @(39, 20, 573.486980432, 20, 407.148412742)
)

foreach ($row in $data) {
    $r = $row[0]

    $ws.Range("A$r").NumberFormat = "0"
    $ws.Range("A$r").Value = $row[1]

    $ws.Range("B$r").NumberFormat = "0.00"
    $ws.Range("B$r").Value = $row[2]

    $ws.Range("D$r").Value = $row[3]

    $ws.Range("E$r").NumberFormat = "0.00"
    $ws.Range("E$r").Value = $row[4]
}

# --- Row 41: "overall" averages for each block (third block is empty -> #DIV/0!) ---
$ws.Range("B41").NumberFormat = "0.00"
$ws.Range("B41").Formula = "=AVERAGE(B20:B39)"

$ws.Range("D41").NumberFormat = "0"
$ws.Range("D41").Value = "overall"
$ws.Range("E41").NumberFormat = "0.00"
$ws.Range("E41").Formula = "=AVERAGE(E20:E39)"

$ws.Range("G41").NumberFormat = "0"
$ws.Range("G41").Value = "overall"
$ws.Range("H41").NumberFormat = "0.00"
$ws.Range("H41").Formula = "=AVERAGE(H20:H39)"

# --- Selection follows where the new data starts, matching the saved view ---
$ws.Range("G20").Select()
